# Adding more info to README and changing organization to Setup guide
#
# Appends four new list-paragraphs to the end of the document body
# (right after the "Java" bullet, before the section properties),
# mirroring the target diff:
#   1. ilvl=1  "This can be done because ps1 script files can be converted into exe files"
#   2. ilvl=0  "Test whether the current configurations are standalone as in they meet these criteria:"
#   3. ilvl=1  "Do not require vcpkg to be installed on the computer and can be run minimally with the
#               asio header file, hidapi header files, and hidapi.dll"
#   4. ilvl=0  (empty paragraph)

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Insert-ListParagraph([string]$innerXml) {
    # Collapsed range right at the very end of the body (before the final
    # paragraph mark) so InsertXML appends a brand-new paragraph instead of
    # replacing existing content.
    $endPos = $d.Content.End - 1
    $rng = $d.Range($endPos, $endPos)
    $rng.InsertXML("<w:p $wNs>$innerXml</w:p>")
}

# 1) ilvl=1 bullet about ps1 -> exe conversion
Insert-ListParagraph '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>This can be done because ps1 script files can be converted into exe files</w:t></w:r>'

# 2) ilvl=0 bullet, with a gramStart/gramEnd proofing mark around "as in"
Insert-ListParagraph '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Test whether the current configurations are standalone </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>as in</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> they meet these criteria:</w:t></w:r>'

# 3) ilvl=1 bullet, with a spellStart/spellEnd proofing mark around "hidapi"
Insert-ListParagraph '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Do not require vcpkg to be installed on the computer and can be run minimally with the asio header file, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hidapi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> header files, and hidapi.dll</w:t></w:r>'

# 4) ilvl=0 empty bullet (new "Setup guide" section placeholder)
Insert-ListParagraph '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'
